$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 310.1
$ws.Range("I33").Value = 200.125
$ws.Range("J33").Value = 750
$ws.Range("K33").Value = 200.125
$ws.Range("L33").Value = 750
$ws.Range("M33").Value = 28.875
$ws.Range("N33").Value = -1208
$ws.Range("H106").Value = 3276
$ws.Range("I106").Value = 2990
$ws.Range("J106").Value = 3466.6667
$ws.Range("K106").Value = 2990
$ws.Range("L106").Value = 3466.6667
$ws.Range("M106").Value = -2359
$ws.Range("N106").Value = -4728.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2578
$ws.Range("I45").Value = 1398
$ws.Range("K45").Value = 1398
$ws.Range("M45").Value = -1021
$ws.Range("H61").Value = 1875
$ws.Range("I61").Value = 1430
$ws.Range("J61").Value = 4100
$ws.Range("K61").Value = 1430
$ws.Range("L61").Value = 4100
$ws.Range("M61").Value = -1218
$ws.Range("N61").Value = -4524
$ws.Range("H74").Value = 2953.2632
$ws.Range("I74").Value = 2317.923
$ws.Range("J74").Value = 4329.8335
$ws.Range("K74").Value = 2317.923
$ws.Range("L74").Value = 4329.8335
$ws.Range("M74").Value = -1443.923
$ws.Range("N74").Value = -6077.8335
$ws.Range("H77").Value = 2953.2632
$ws.Range("I77").Value = 2317.923
$ws.Range("J77").Value = 4329.8335
$ws.Range("K77").Value = 11589.615
$ws.Range("L77").Value = 21649.1675
$ws.Range("M77").Value = -7221.614999999998
$ws.Range("N77").Value = -30385.1675
$ws.Range("H110").Value = 2573.9092
$ws.Range("I110").Value = 2573.9092
$ws.Range("K110").Value = 2573.9092
$ws.Range("M110").Value = -528.9092000000001
$ws.Range("H122").Value = 3127.6365
$ws.Range("J122").Value = 9000
$ws.Range("L122").Value = 27000
$ws.Range("N122").Value = -31900
$ws.Range("H132").Value = 2664.3572
$ws.Range("I132").Value = 1267.9445
$ws.Range("J132").Value = 5177.9
$ws.Range("K132").Value = 3803.8335
$ws.Range("L132").Value = 15533.7
$ws.Range("M132").Value = -1273.8335
$ws.Range("N132").Value = -20593.7
$ws.Range("H136").Value = 1875
$ws.Range("I136").Value = 1430
$ws.Range("J136").Value = 4100
$ws.Range("K136").Value = 4290
$ws.Range("L136").Value = 12300
$ws.Range("M136").Value = -1740
$ws.Range("N136").Value = -17400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1047.75
$ws.Range("I107").Value = 1047.75
$ws.Range("K107").Value = 1047.75
$ws.Range("M107").Value = 872.25
$ws.Range("H134").Value = 3007.9614
$ws.Range("I134").Value = 1758.4722
$ws.Range("J134").Value = 5819.3125
$ws.Range("K134").Value = 5275.4166
$ws.Range("L134").Value = 17457.9375
$ws.Range("M134").Value = -2740.4166
$ws.Range("N134").Value = -22527.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2073.5574
$ws.Range("I58").Value = 1715.2931
$ws.Range("J58").Value = 9000
$ws.Range("K58").Value = 1715.2931
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -1512.2931
$ws.Range("N58").Value = -9406
$ws.Range("H132").Value = 2351.5117
$ws.Range("I132").Value = 1129.2812
$ws.Range("J132").Value = 5907.091
$ws.Range("K132").Value = 3387.8436
$ws.Range("L132").Value = 17721.273
$ws.Range("M132").Value = -857.8435999999997
$ws.Range("N132").Value = -22781.273
$ws.Range("H134").Value = 9106.5
$ws.Range("I134").Value = 12348.333
$ws.Range("J134").Value = 4938.4287
$ws.Range("K134").Value = 37044.999
$ws.Range("L134").Value = 14815.2861
$ws.Range("M134").Value = -34509.999
$ws.Range("N134").Value = -19885.2861
$ws.Range("H136").Value = 2073.5574
$ws.Range("I136").Value = 1715.2931
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 5145.879300000001
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -2595.879300000001
$ws.Range("N136").Value = -32100
$ws.Range("H141").Value = 15908.823
$ws.Range("J141").Value = 15908.823
$ws.Range("L141").Value = 15908.823
$ws.Range("N141").Value = -26268.823

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2409.5186
$ws.Range("I140").Value = 1959.8182
$ws.Range("K140").Value = 5879.4546
$ws.Range("M140").Value = -699.4546
$ws.Range("H141").Value = 7173.278
$ws.Range("I141").Value = 6201.9
$ws.Range("K141").Value = 18605.7
$ws.Range("M141").Value = -13425.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 31252876
$ws.Range("J80").Value = 3286.5715
$ws.Range("L80").Value = 3286.5715
$ws.Range("N80").Value = -5282.5715
$ws.Range("H83").Value = 31252876
$ws.Range("J83").Value = 3286.5715
$ws.Range("L83").Value = 16432.8575
$ws.Range("N83").Value = -26416.8575
$ws.Range("H132").Value = 4997.357
$ws.Range("I132").Value = 4674.6665
$ws.Range("K132").Value = 14023.9995
$ws.Range("M132").Value = -11493.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5383.4375
$ws.Range("I7").Value = 3785.7144
$ws.Range("J7").Value = 6626.1113
$ws.Range("K7").Value = 3785.7144
$ws.Range("L7").Value = 6626.1113
$ws.Range("M7").Value = -3673.7144
$ws.Range("N7").Value = -6850.1113
$ws.Range("H68").Value = 905.9318
$ws.Range("I68").Value = 831.55
$ws.Range("J68").Value = 1649.75
$ws.Range("K68").Value = 831.55
$ws.Range("L68").Value = 1649.75
$ws.Range("M68").Value = -82.54999999999995
$ws.Range("N68").Value = -3147.75
$ws.Range("H71").Value = 905.9318
$ws.Range("I71").Value = 831.55
$ws.Range("J71").Value = 1649.75
$ws.Range("K71").Value = 4157.75
$ws.Range("L71").Value = 8248.75
$ws.Range("M71").Value = -413.75
$ws.Range("N71").Value = -15736.75
$ws.Range("H100").Value = 2184.2666
$ws.Range("I100").Value = 2028.8889
$ws.Range("J100").Value = 2417.3333
$ws.Range("K100").Value = 2028.8889
$ws.Range("L100").Value = 2417.3333
$ws.Range("M100").Value = -1487.8889
$ws.Range("N100").Value = -3499.3333
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 5383.4375
$ws.Range("I126").Value = 3785.7144
$ws.Range("J126").Value = 6626.1113
$ws.Range("K126").Value = 11357.1432
$ws.Range("L126").Value = 19878.3339
$ws.Range("M126").Value = -8887.143199999999
$ws.Range("N126").Value = -24818.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6292727
$ws.Range("I132").Value = 3819.9033
$ws.Range("J132").Value = 15154369
$ws.Range("K132").Value = 11459.7099
$ws.Range("L132").Value = 45463107
$ws.Range("M132").Value = -8929.7099
$ws.Range("N132").Value = -45468167
$ws.Range("H136").Value = 6378.067
$ws.Range("I136").Value = 1174.75
$ws.Range("J136").Value = 8270.182000000001
$ws.Range("K136").Value = 3524.25
$ws.Range("L136").Value = 24810.546
$ws.Range("M136").Value = -974.25
$ws.Range("N136").Value = -29910.546
